# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Intel(R) Dual Band Wireless-AC 8265 - 20.70.27.1
$ws.Range("C3").Value = 702

# Row 4: Intel(R) Wi-Fi 6 AX200 160MHz - 23.80.0.7
$ws.Range("C4").Value = 1340

# Row 5: Intel(R) Dual Band Wireless-AC 8265 - 20.70.18.2
$ws.Range("C5").Value = 1850

# Row 6: Intel(R) Dual Band Wireless-AC 8265 - 20.70.30.1
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = 1104

# Row 7: Totals
$ws.Range("B7").Value = 78
$ws.Range("C7").Value = 4996
